# Fix "Closing Date" value to include a missing space, and fill in
# "N/A" for the previously-empty "Site Briefing Date" and "Remarks" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Closing Date value (row 6) was missing a space between the date and time.
$ws.Range("B6").Value = "07 Apr 2022 01:00PM"

# Site Briefing Date (row 7) and Remarks (row 8) were blank; set to N/A.
$ws.Range("B7").Value = "N/A"
$ws.Range("B8").Value = "N/A"
